$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$origStyle = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.977.15"
$ws.Range("D2").Style = $origStyle
$ws.Range("E2").Value = "  -1.07%  "
$origStyle = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.365.51"
$ws.Range("D3").Style = $origStyle
$ws.Range("E3").Value = "  -1.12%  "
$ws.Range("E4").Value = "  +0.17%  "
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "405.61"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  -1.39%  "
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.61"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  +9.85%  "
$ws.Range("E7").Value = "  +2.39%  "
$ws.Range("E8").Value = "  +0.09%  "
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.671"
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = "  +4.82%  "
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.121"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  +7.32%  "
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.38"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  +2.45%  "
$ws.Range("E12").Value = "  -0.76%  "
$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.890.51"
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = "  -1.22%  "
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.31"
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = "  -0.92%  "
$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "19.57"
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = "  +0.52%  "
$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.367.42"
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = "  -1.85%  "
$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.013.20"
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = "  -0.75%  "
$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.01"
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = "  -0.55%  "
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.04"
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = "  +1.64%  "
$ws.Range("E20").Value = "  +8.03%  "
$ws.Range("E21").Value = "  -3.39%  "
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "84.23"
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = "  +10.71%  "
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "308.68"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  +2.65%  "
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.68"
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = "  -1.05%  "
$ws.Range("E25").Value = "  +0.19%  "
$ws.Range("E26").Value = "  +11.89%  "
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.27"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  +8.05%  "
$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "29.34"
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = "  -4.62%  "
$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.53"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  -7.73%  "
$ws.Range("E30").Value = "  +1.00%  "
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.117"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  +1.48%  "
$ws.Range("E32").Value = "  +0.00%  "
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.26"
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = "  -1.30%  "
$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "41.09"
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = "  -3.02%  "
$ws.Range("E35").Value = "  -1.03%  "
$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0479"
$ws.Range("D36").Style = $origStyle
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.86"
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = "  -0.89%  "
$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.998"
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = "  +0.11%  "
$ws.Range("E39").Value = "  -2.73%  "
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.89"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  -4.35%  "
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.98"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  +0.86%  "
$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "136.79"
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = "  +2.13%  "
$ws.Range("E43").Value = "  +0.95%  "
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.03"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  +2.97%  "
$ws.Range("E45").Value = "  +0.84%  "
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.61"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  -4.01%  "
$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.22"
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = "  +1.35%  "
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.31"
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = "  -2.61%  "
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.119.90"
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = "  -3.94%  "
$ws.Range("E50").Value = "  -3.84%  "
$ws.Range("E51").Value = "  +0.02%  "
